$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.647.34'
$ws.Range('E2').Value = '  -1.44%  '
$ws.Range('D3').Value = '2.443.97'
$ws.Range('E3').Value = '  -1.36%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = "'569.76"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('D6').Value = "'144.20"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.85%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').Value = "'0.533"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.47%  '
$ws.Range('D9').Value = '2.444.14'
$ws.Range('E9').Value = '  -1.78%  '
$ws.Range('D10').Value = "'0.109"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.88%  '
$ws.Range('E11').Value = '  +1.49%  '
$ws.Range('D12').Value = "'5.24"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.75%  '
$ws.Range('D13').Value = "'0.355"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.59%  '
$ws.Range('D14').Value = "'27.02"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').Value = "'0.0000175"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.94%  '
$ws.Range('D16').Value = '2.889.70'
$ws.Range('E16').Value = '  -0.68%  '
$ws.Range('D17').Value = '62.479.36'
$ws.Range('E17').Value = '  -1.37%  '
$ws.Range('D18').Value = '2.449.25'
$ws.Range('E18').Value = '  -1.54%  '
$ws.Range('D19').Value = "'11.22"
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').Value = "'7.27"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.32%  '
$ws.Range('D21').Value = "'327.16"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.58%  '
$ws.Range('D22').Value = "'4.17"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.87%  '
$ws.Range('D23').Value = "'2.11"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +10.43%  '
$ws.Range('E24').Value = '  +0.48%  '
$ws.Range('D25').Value = "'65.33"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.42%  '
$ws.Range('D26').Value = "'623.34"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.39%  '
$ws.Range('D27').Value = "'9.04"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.12%  '
$ws.Range('D28').Value = '0.0₃0997'
$ws.Range('E28').Value = '  -5.18%  '
$ws.Range('E29').Value = '  -2.69%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = "'0.999"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.83%  '
$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').Value = "'1.49"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.49%  '
$ws.Range('D32').Value = "'8.14"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.49%  '
$ws.Range('E33').Value = '  -1.58%  '
$ws.Range('E34').Value = '  -4.20%  '
$ws.Range('D35').Value = "'5.12"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.11%  '
$ws.Range('D36').Value = "'1.50"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.07%  '
$ws.Range('E37').Value = '  +0.32%  '
$ws.Range('D38').Value = "'0.377"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.69%  '
$ws.Range('D39').Value = "'18.82"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.81%  '
$ws.Range('D40').Value = "'5.34"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.04%  '
$ws.Range('D41').Value = "'146.44"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.30%  '
$ws.Range('D42').Value = "'1.76"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.76%  '
$ws.Range('D43').Value = "'2.57"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.28%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = "'42.22"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.05%  '
$ws.Range('B45').Value = 'USDe'
$ws.Range('C45').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D45').Value = "'0.999"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.01%  '
$ws.Range('D46').Value = "'146.70"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.30%  '
$ws.Range('D47').Value = "'3.76"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.03%  '
$ws.Range('D48').Value = "'20.71"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.90%  '
$ws.Range('E49').Value = '  -4.13%  '
$ws.Range('E50').Value = '  -1.98%  '
$ws.Range('D51').Value = "'0.0231"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.96%  '
